$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.2159279182453454
$ws.Range("D2").Value = 0.831033411157537

# Row 3
$ws.Range("C3").Value = -1.025707438976052
$ws.Range("D3").Value = 0.3161760420499813

# Row 4
$ws.Range("C4").Value = -2.743840755384584
$ws.Range("D4").Value = 0.0118500545612259

# Row 5
$ws.Range("C5").Value = -4.669040918040095
$ws.Range("D5").Value = 0.0001177545371808808

# Row 6
$ws.Range("C6").Value = -0.7377658936177198
$ws.Range("D6").Value = 0.4684536135066821

# Row 7
$ws.Range("C7").Value = -2.121469381947279
$ws.Range("D7").Value = 0.04538612117626628

# Row 8
$ws.Range("C8").Value = -4.10053781979866
$ws.Range("D8").Value = 0.0004719726044126737

# Row 9
$ws.Range("C9").Value = -1.733804290229447
$ws.Range("D9").Value = 0.09694659853638954
$ws.Range("G9").Value = "No"

# Row 10
$ws.Range("C10").Value = -3.055892475233526
$ws.Range("D10").Value = 0.005791538328698831

# Row 11
$ws.Range("C11").Value = -0.8832288212678392
$ws.Range("D11").Value = 0.3866625721520309
